# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计", before "2022-Q2".
# 2) Populate it with the single fund row of new quarterly data.
# 3) Insert a new top row in "总计" summarizing the 2022-Q3 quarter, pushing
#    the existing rows down by one, and renumber the index column (A).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$wb.Worksheets.Add($null, $wb.Worksheets.Item("总计")) | Out-Null
$wb.Worksheets.Item(2).Name = "2022-Q3"

# Copy header formatting (bordered/bold/centered style) from the
# equivalent quarterly sheet so the new sheet matches the others.
# (Sheet references are re-fetched by name after the Add() above, since
# handles captured before a sheet-collection change go stale.)
$wb.Worksheets.Item("2022-Q2").Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("2022-Q2").Range("A2").Copy()
$wb.Worksheets.Item("2022-Q3").Range("A2").PasteSpecial(-4122)

$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data row - text-like numeric columns (fund code / amounts) are stored
# as text, matching the other quarterly sheets; only the rank column is
# a real number.
$q3Sheet.Cells.Item(2, 1).Value = 0

$q3Sheet.Cells.Item(2, 2).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 2).Value = "162102"
$q3Sheet.Cells.Item(2, 2).Style = "Normal"

$q3Sheet.Cells.Item(2, 3).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 3).Value = "金鹰中小盘精选混合"
$q3Sheet.Cells.Item(2, 3).Style = "Normal"

$q3Sheet.Cells.Item(2, 4).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 4).Value = "3.17"
$q3Sheet.Cells.Item(2, 4).Style = "Normal"

$q3Sheet.Cells.Item(2, 5).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 5).Value = "78.28"
$q3Sheet.Cells.Item(2, 5).Style = "Normal"

$q3Sheet.Cells.Item(2, 6).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 6).Value = "4.30"
$q3Sheet.Cells.Item(2, 6).Style = "Normal"

$q3Sheet.Cells.Item(2, 7).NumberFormat = "@"
$q3Sheet.Cells.Item(2, 7).Value = "0.1363"
$q3Sheet.Cells.Item(2, 7).Style = "Normal"

$q3Sheet.Cells.Item(2, 8).Value = 10

# ---------------------------------------------------------------------
# Step 2: insert the new summary row into "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.14

# Match column-A's bordered/bold index style on the new row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the zero-based index column for the rows pushed down.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "2022-Q3 data added"
